$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptos list.
$ws.Range("D2").Value = '70.636.69'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Value = '2.525.73'
$ws.Range("E3").Value = '  -5.31%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''574.90'
$ws.Range("E5").Value = '  -4.08%  '
$ws.Range("D6").Value = '''169.71'
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''0.508'
$ws.Range("D9").Value = '2.525.83'
$ws.Range("E9").Value = '  -5.28%  '
$ws.Range("E10").Value = '  -2.83%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E13").Value = '  -3.77%  '
$ws.Range("D14").Value = '2.988.17'
$ws.Range("E14").Value = '  -5.79%  '
$ws.Range("D15").Value = '70.511.52'
$ws.Range("E15").Value = '  -1.90%  '
$ws.Range("E16").Value = '  -2.84%  '
$ws.Range("D17").Value = '''24.85'
$ws.Range("E17").Value = '  -5.23%  '
$ws.Range("D18").Value = '2.527.44'
$ws.Range("E18").Value = '  -5.62%  '
$ws.Range("D19").Value = '''11.53'
$ws.Range("E19").Value = '  -5.62%  '
$ws.Range("D20").Value = '''7.53'
$ws.Range("E20").Value = '  -8.31%  '
$ws.Range("D21").Value = '''356.46'
$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("E22").Value = '  -5.92%  '
$ws.Range("E23").Value = '  -3.89%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '''69.26'
$ws.Range("E25").Value = '  -3.74%  '
$ws.Range("E26").Value = '  -6.52%  '
$ws.Range("D27").Value = '''9.21'
$ws.Range("E27").Value = '  -5.92%  '
$ws.Range("E28").Value = '  -5.68%  '
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("E30").Value = '  -6.39%  '
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").Value = '''479.84'
$ws.Range("E32").Value = '  -4.21%  '
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("D34").Value = '''1.75'
$ws.Range("E34").Value = '  -3.84%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '''157.68'
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("E37").Value = '  +4.87%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '''18.58'
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = '''18.84'
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  -5.07%  '
$ws.Range("D42").Value = '''1.65'
$ws.Range("E42").Value = '  -6.93%  '
$ws.Range("E43").Value = '  -4.36%  '
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("E45").Value = '  -6.12%  '
$ws.Range("D46").Value = '''38.32'
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("D47").Value = '''143.04'
$ws.Range("E47").Value = '  -7.83%  '
$ws.Range("E48").Value = '  -5.42%  '
$ws.Range("E49").Value = '  -6.63%  '
$ws.Range("E50").Value = '  -6.77%  '
$ws.Range("E51").Value = '  -1.51%  '

# The apostrophe trick above stamps a quote-prefix format on each cell.
# Reset those specific cells back to the default "Normal" style so they
# match the plain, unstyled cells used throughout the rest of the sheet.
$textForcedCells = @("D4", "D5", "D6", "D8", "D17", "D19", "D20", "D21", "D25", "D27", "D32", "D34", "D35", "D36", "D38", "D39", "D42", "D46", "D47")
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
